# Ran model for 1/14/2021
#
# - Row 65 (PHO @ ATL) turned out to be a bad/duplicate entry from the
#   previous model run; remove it, which shifts the remaining 1/13/2021
#   rows up by one (old rows 66-67 become 65-66).
# - Fill in the "Beat Vegas?" (col G) verdict for every 1/13/2021 game now
#   that final scores are known.
# - Append the new model predictions for the 1/14/2021 slate as rows 67-71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stale PHO @ ATL row -----------------------------------
$ws.Rows.Item(65).Delete()

# --- Beat Vegas? verdicts for 1/13/2021 games (rows 60-66 post-shift) --
$ws.Range("G60").Value = "No"
$ws.Range("G61").Value = "No"
$ws.Range("G62").Value = "Yes"
$ws.Range("G63").Value = "Yes"
$ws.Range("G64").Value = "No"
$ws.Range("G65").Value = "Yes"
$ws.Range("G66").Value = "No"

# --- New rows for 1/14/2021, matching the date column's formatting ----
$ws.Range("A66").Copy()
$ws.Range("A67:A71").PasteSpecial(-4122)

$ws.Range("A67").Value = 44210
$ws.Range("B67").Value = "PHI"
$ws.Range("C67").Value = "MIA"
$ws.Range("D67").Value = -7
$ws.Range("E67").Value = 0.6
$ws.Range("F67").Value = -7.6

$ws.Range("A68").Value = 44210
$ws.Range("B68").Value = "SAS"
$ws.Range("C68").Value = "HOU"
$ws.Range("D68").Value = -7
$ws.Range("E68").Value = 2.6
$ws.Range("F68").Value = -9.6

$ws.Range("A69").Value = 44210
$ws.Range("B69").Value = "CHO"
$ws.Range("C69").Value = "TOR"
$ws.Range("D69").Value = 10
$ws.Range("E69").Value = 4.3
$ws.Range("F69").Value = 5.7

$ws.Range("A70").Value = 44210
$ws.Range("B70").Value = "DEN"
$ws.Range("C70").Value = "GSW"
$ws.Range("D70").Value = -4.5
$ws.Range("E70").Value = -3.4
$ws.Range("F70").Value = -1.1

$ws.Range("A71").Value = 44210
$ws.Range("B71").Value = "POR"
$ws.Range("C71").Value = "IND"
$ws.Range("D71").Value = -3
$ws.Range("E71").Value = -3.3
$ws.Range("F71").Value = 0.3

# --- Restore the view: scrolled to row 53, cursor left on I71 ---------
$excel.ActiveWindow.ScrollRow = 53
$ws.Range("I71").Select()
